# Doing Updates for Financials
# Insert a new "most recent period" column before column D on the HWCC sheet,
# shifting the existing D:K data right to E:L, then populate the new column D
# with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HWCC")

# Insert a new column before D; existing D:K (and their formatting) shift to E:L.
$ws.Columns("D").Insert()

# Copy number formatting from the (now shifted) column E into the new column D
# so the new column matches the date/number styling used throughout the sheet.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Income Statement (new period column D) ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 356900
$ws.Range("D9").Value = 271700
$ws.Range("D10").Value = 85200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 100
$ws.Range("D15").Value = 2200
$ws.Range("D17").Value = 343000
$ws.Range("D18").Value = 13900
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 16100
$ws.Range("D22").Value = 2900
$ws.Range("D23").Value = 11000
$ws.Range("D24").Value = 2400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 8600
$ws.Range("D27").Value = 8600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 8600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 8600

# --- Balance Sheet (new period column D) ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 1400
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 60200
$ws.Range("D44").Value = 94300
$ws.Range("D45").Value = 700
$ws.Range("D46").Value = 156700
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 11500
$ws.Range("D49").Value = 33500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 203100
$ws.Range("D57").Value = 11300
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 19200
$ws.Range("D60").Value = 30500
$ws.Range("D61").Value = 71300
$ws.Range("D62").Value = 600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 102400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 106000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 100700
$ws.Range("D77").Value = 0

# --- Cash Flow Statement (new period column D) ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 8600
$ws.Range("D83").Value = 2200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 5300
$ws.Range("D91").Value = -1500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1500
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -2500
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 1400
